# The deck currently applies the "Integral" design/theme to its slide
# master (and therefore to every slide). This edit switches the
# presentation's theme colors back to the default "Office Theme" palette
# (as if "Office Theme" had been picked from the Design tab's theme
# gallery), by rewriting the master's 12 theme colors in place.
#
# PowerPoint's ColorScheme.Colors(index) follows the
# MsoThemeColorSchemeIndex order, matching <a:clrScheme> 1:1:
#   1 dk1   2 lt1   3 dk2     4 lt2     5 accent1  6 accent2
#   7 accent3  8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink

$p = $ppt.ActivePresentation

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme palette (the built-in default PowerPoint theme).
$officeThemeColors = @(
    (RGBVal 0x00 0x00 0x00),  # dk1
    (RGBVal 0xFF 0xFF 0xFF),  # lt1
    (RGBVal 0x44 0x54 0x6A),  # dk2
    (RGBVal 0xE7 0xE6 0xE6),  # lt2
    (RGBVal 0x5B 0x9B 0xD5),  # accent1
    (RGBVal 0xED 0x7D 0x31),  # accent2
    (RGBVal 0xA5 0xA5 0xA5),  # accent3
    (RGBVal 0xFF 0xC0 0x00),  # accent4
    (RGBVal 0x44 0x72 0xC4),  # accent5
    (RGBVal 0x70 0xAD 0x47),  # accent6
    (RGBVal 0x05 0x63 0xC1),  # hlink
    (RGBVal 0x95 0x4F 0x72)   # folHlink
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
